$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking text (e.g. "228.99") need the
# column pre-set to Text format ("@") before assignment, otherwise Excel
# auto-converts the input to a real number (losing the inlineStr/text type
# and introducing float rounding). Apply "@" first to each such cell (the
# identical format string de-dupes to a single shared style), then write
# the text value.
$textCells = @(
    "D5"
    "D6"
    "D7"
    "D10"
    "D12"
    "D14"
    "D19"
    "D20"
    "D22"
    "D25"
    "D26"
    "D27"
    "D29"
    "D30"
    "D33"
    "D34"
    "D36"
    "D38"
    "D40"
    "D41"
    "D42"
    "D49"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- value updates ---
$ws.Range("D2").Value = "39.396.94"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "2.158.45"
$ws.Range("E3").Value = "  +3.50%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "228.99"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "0.621"
$ws.Range("D7").Value = "63.05"
$ws.Range("E7").Value = "  +4.65%  "
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("D10").Value = "0.0866"
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "15.99"
$ws.Range("E12").Value = "  +6.82%  "
$ws.Range("D13").Value = "2.479.43"
$ws.Range("E13").Value = "  +3.51%  "
$ws.Range("D14").Value = "22.15"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "2.158.53"
$ws.Range("E17").Value = "  +3.27%  "
$ws.Range("D18").Value = "39.357.73"
$ws.Range("D19").Value = "72.28"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").Value = "6.13"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").Value = "228.72"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "9.62"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("D27").Value = "171.84"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "19.67"
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "1.42"
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("E31").Value = "  +8.61%  "
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").Value = "4.66"
$ws.Range("E33").Value = "  +3.57%  "
$ws.Range("D34").Value = "4.81"
$ws.Range("E34").Value = "  +2.73%  "
$ws.Range("E35").Value = "  +9.41%  "
$ws.Range("D36").Value = "0.0621"
$ws.Range("E36").Value = "  +2.12%  "
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("D38").Value = "3.57"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "103.81"
$ws.Range("E40").Value = "  +3.19%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0231"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "18.06"
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("D43").Value = "1.533.85"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").Value = "  +5.95%  "
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("E46").Value = "  +7.06%  "
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("D49").Value = "4.21"
$ws.Range("E49").Value = "  +2.51%  "
$ws.Range("D50").Value = "2.363.37"
$ws.Range("E50").Value = "  +3.46%  "
$ws.Range("E51").Value = "  +0.21%  "
